$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44754
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 16000
$ws.Range("L2").Value = 17000
$ws.Range("M2").Value = 16500
$ws.Range("P2").Value = 917

$ws.Range("D3").Value = 44637
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15500
$ws.Range("P3").Value = 861

$ws.Range("D4").Value = 44630
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 861

$ws.Range("D5").Value = 44635
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 15500
$ws.Range("P5").Value = 861

$ws.Range("D6").Value = 44628
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("P6").Value = 861

$ws.Range("D7").Value = 44645
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 16000
$ws.Range("M7").Value = 15500
$ws.Range("P7").Value = 861

$ws.Range("D8").Value = 44769
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 17000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 17500
$ws.Range("P8").Value = 972

$ws.Range("D9").Value = 44642
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("P9").Value = 861

$ws.Range("D10").Value = 44659
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 15000
$ws.Range("L10").Value = 16000
$ws.Range("M10").Value = 15500
$ws.Range("P10").Value = 861

$ws.Range("D11").Value = 44649
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15500
$ws.Range("P11").Value = 861

$ws.Range("D12").Value = 44658
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 15000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 15500
$ws.Range("P12").Value = 861

$ws.Range("D13").Value = 44785
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 17500
$ws.Range("P13").Value = 972

$ws.Range("D14").Value = 44384
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 17000
$ws.Range("L14").Value = 18000
$ws.Range("M14").Value = 17500
$ws.Range("P14").Value = 972

$ws.Range("D15").Value = 44384
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 15000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 15000
$ws.Range("P15").Value = 833

$ws.Range("D16").Value = 44664
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 16000
$ws.Range("M16").Value = 15500
$ws.Range("P16").Value = 861

$ws.Range("D17").Value = 44651
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("P17").Value = 861

$ws.Range("D18").Value = 44656
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15500
$ws.Range("P18").Value = 861

$ws.Range("D19").Value = 44775
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 17500
$ws.Range("P19").Value = 972

$ws.Range("D20").Value = 44761
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17500
$ws.Range("P20").Value = 972

$ws.Range("D21").Value = 44763
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("P21").Value = 972

$ws.Range("D22").Value = 44771
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17500
$ws.Range("P22").Value = 972

$ws.Range("D23").Value = 44782
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 17000
$ws.Range("L23").Value = 18000
$ws.Range("M23").Value = 17500
$ws.Range("P23").Value = 972
